# Daily attendance processing - 2025-12-31 09:59:24
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Top "Class Statistics" summary box (K/L columns)
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 201        # Recorded Sessions
$ws.Range("L7").Value = 3          # Missing Sessions
$ws.Range("L9").Value = "'63.2%"   # Coverage %
$ws.Range("L10").Value = "'75.8%"  # Average Attendance %

# ---------------------------------------------------------------------------
# 2. "Group Statistics" breakdown table rows 21-26 (O/P/R/S columns)
# ---------------------------------------------------------------------------
$groupStats = @(
    @{ Row = 21; O = 17; P = 0; R = "63.0%"; S = "78.3%" },
    @{ Row = 22; O = 17; P = 0; R = "63.0%"; S = "77.1%" },
    @{ Row = 23; O = 17; P = 0; R = "63.0%"; S = "79.9%" },
    @{ Row = 24; O = 16; P = 1; R = "59.3%"; S = "72.5%" },
    @{ Row = 25; O = 17; P = 0; R = "63.0%"; S = "70.8%" },
    @{ Row = 26; O = 17; P = 0; R = "63.0%"; S = "63.9%" }
)

foreach ($stat in $groupStats) {
    $r = $stat.Row
    $ws.Range("O$r").Value = $stat.O
    $ws.Range("P$r").Value = $stat.P
    $ws.Range("R$r").Value = "'" + $stat.R
    $ws.Range("S$r").Value = "'" + $stat.S
}

# ---------------------------------------------------------------------------
# 3. Swap "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
#    in the "Recorded By" (G) column across the sheet.
# ---------------------------------------------------------------------------
$recordedByRows = @(8,9,10,12,14,15,17,18,34,35,36,38,40,41,43,44,60,61,62,64,
    66,67,69,70,86,87,88,90,92,93,95,96,112,113,114,116,118,119,121,122,138,
    139,140,142,144,145,147,148,164,167,170,191,194,197,218,221,224,245,248,
    251,272,275,278,299,302,305)

foreach ($r in $recordedByRows) {
    $ws.Range("G$r").Value = "dnasr281@gmail.com, System"
}

# ---------------------------------------------------------------------------
# 4. Newly-recorded sessions (row 17, 31/12/2025) for groups B1D1, B1D2,
#    B1E1, B1E2, B1F1, B1F2 -- were "Not Recorded" (pink), now "Recorded"
#    (green), with a recorder e-mail and real attendance counts.
# ---------------------------------------------------------------------------
$newlyRecorded = @(
    @{ Row = 174; Present = 19; Total = 23 },
    @{ Row = 201; Present = 21; Total = 30 },
    @{ Row = 228; Present = 19; Total = 26 },
    @{ Row = 255; Present = 24; Total = 28 },
    @{ Row = 282; Present = 20; Total = 26 },
    @{ Row = 309; Present = 25; Total = 29 }
)

foreach ($rec in $newlyRecorded) {
    $r = $rec.Row
    $rowRange = $ws.Range("A" + $r + ":I" + $r)
    $rowRange.Interior.Color = 9498256   # green fill, same as other "Recorded" rows

    $ws.Range("G$r").Value = "dnasr281@gmail.com"
    $ws.Range("H$r").Value = "$($rec.Present)/$($rec.Total)"
    $ws.Range("I$r").Value = "Recorded"
}
